$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the requirement descriptions (column D) to reflect the
#     "usuario-cliente" -> "usuario-empresa" terminology change and the
#     corrected "UC.002 Manter Transportadora" text.
$ws.Range("D3").Value = "O sistema deve permitir que o usuário-empresa cadastre seus dados, sendo estes, CNPJ (não modificável), Razão Social, Endereço, Número de endereço, Setor, Cidade, UF, CEP, Fone, Contato e E-mail de contato. Ao efetuar o cadastro, o usuário-empresa deve assinar o contrato via online e o sistema deve validar se o CNPJ é existente."

$ws.Range("D4").Value = "O sistema deve permitir que o usuário-transportador cadastre e mantenha seus dados, sendo estes, CNPJ (não modificável), Razão Social, Nome Fantasia, Inscrição Estadual, Endereço, Número de endereço, Setor, Cidade, UF, CEP, Fone, Contato, E-mail de contato. Ao efetuar o cadastro, o usuário-transportador deve assinar o contrato via online e o sistema deve validar se o CNPJ é existente."

$ws.Range("D5").Value = "O sistema deve permitir que o usuário-empresa cadastre  remessas para cotações, devendo conter as seguintes informações: Tipo de mercadoria*, Peso Bruto, Quantidade de volumes ou de pallets, Valor Total da Mercadoria, Dimensões dos volumes, Endereço, Cidade Origem*, Prazo mínimo para pagamento de boleto, Nome da pessoa para contato, Telefone de contato, Dia da disponibilidade de coleta, Tempo de duração para receber lances*. O usuário-empresa deve optar por ser Emitente ou Destinatário. Caso o usuário-empresa não queira preencher todos os campos, poderá estar disponibilizando as informações através de anexo ou colar no campo Informações em formato texto. Todas as informações com * deverá ser preenchida obrigatoriamente."

$ws.Range("C6").Value = "UC.004 Emitir Lance"
$ws.Range("D6").Value = "O sistema deve permitir que o usuário-transportador emite lance em remessas postadas por usuário-empresa."

$ws.Range("C7").Value = "UC.005 Visualizar Lances"
$ws.Range("D7").Value = "O sistema deve permitir que o usuário-empresa visualize todas os lances recebidos em sua remessa."

$ws.Range("C8").Value = "UC.007 Iniciar Contrato"
$ws.Range("D8").Value = "Este requisito tem como objetivo permitir que o usuário empresa inicie um contrato de frete a partir de um lance recebido em sua remessa."

# --- Adjust row heights for the rows whose text got shorter.
$ws.Rows(4).RowHeight = 89.25
$ws.Rows(5).RowHeight = 165.75
$ws.Rows(8).RowHeight = 38.25

# --- Remove the old "UC.006 Metodo Pagamento" row (row 9); this also
#     shrinks the table range and shifts blank rows up.
$ws.Rows(9).Delete()

# --- Restore the sheet view: show from the top (A1), zoom to 115%,
#     and select B8:E8.
$ws.Activate()
$ws.Range("B8:E8").Select()
$excel.ActiveWindow.Zoom = 115
